$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.959.78"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").Value = "2.317.45"
$ws.Range("E3").Value = "  +1.72%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.76"
$ws.Range("E5").Value = "  +20.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.52"
$ws.Range("E6").Value = "  +1.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +0.81%  "

$ws.Range("E8").Value = "  +0.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.624"
$ws.Range("E9").Value = "  +2.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.41"
$ws.Range("E10").Value = "  +7.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.89"
$ws.Range("E12").Value = "  +15.30%  "

$ws.Range("E13").Value = "  +2.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.77"
$ws.Range("E14").Value = "  +4.02%  "

$ws.Range("D15").Value = "2.662.55"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.857"
$ws.Range("E16").Value = "  +1.24%  "

$ws.Range("D17").Value = "2.319.95"
$ws.Range("E17").Value = "  +1.33%  "

$ws.Range("D18").Value = "43.873.82"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("E19").Value = "  +2.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.69"
$ws.Range("E20").Value = "  +8.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.80"
$ws.Range("E21").Value = "  +1.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.49"
$ws.Range("E22").Value = "  +5.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.62"
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.89"
$ws.Range("E24").Value = "  +16.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.42"
$ws.Range("E25").Value = "  +5.35%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.48"
$ws.Range("E27").Value = "  +0.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.69"
$ws.Range("E28").Value = "  +10.74%  "

$ws.Range("E29").Value = "  -0.44%  "

$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "177.94"
$ws.Range("E31").Value = "  +0.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.99"
$ws.Range("E32").Value = "  -0.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0936"
$ws.Range("E33").Value = "  +5.72%  "

$ws.Range("E34").Value = "  +4.43%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("E35").Value = "  +0.73%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.74"
$ws.Range("E36").Value = "  +6.76%  "

$ws.Range("E37").Value = "  +3.58%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.95"
$ws.Range("E38").Value = "  +20.62%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0357"
$ws.Range("E39").Value = "  +0.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.245"
$ws.Range("E40").Value = "  +3.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.40"
$ws.Range("E41").Value = "  +1.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.10"
$ws.Range("E42").Value = "  +11.71%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.41"
$ws.Range("E44").Value = "  +3.83%  "

$ws.Range("E45").Value = "  +6.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.93"
$ws.Range("E46").Value = "  +13.82%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.83"
$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("E48").Value = "  -0.97%  "

$ws.Range("B49").Value = "WOONetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.472"
$ws.Range("E49").Value = "  +11.63%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "100.36"
$ws.Range("E50").Value = "  +1.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.23"
$ws.Range("E51").Value = "  +3.09%  "

